$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(290).Insert()
$ws.Range("A290").Value = 3
$ws.Range("B290").Value = "Femacal de La Calera"
$ws.Range("C290").Value = "Coquimbo"
$ws.Range("D290").Value = 44825
$ws.Range("E290").Value = 5
$ws.Range("F290").Value = 100112009
$ws.Range("G290").Value = "Acelga"
$ws.Range("H290").Value = "Sin especificar"
$ws.Range("I290").Value = "Primera"
$ws.Range("J290").Value = 230
$ws.Range("K290").Value = 2800
$ws.Range("L290").Value = 3000
$ws.Range("M290").Value = 2904
$ws.Range("N290").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O290").Value = "Provincia de Quillota"
$ws.Range("P290").Value = 484
$ws.Range("Q290").Value = 6
$ws.Range("R290").Value = "Hortaliza"
